$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values: rows 2-7 become "2022/2010", rows 8-13 become "2022/2021"
$ws.Range("A2:A7").Value = "2022/2010"
$ws.Range("A8:A13").Value = "2022/2021"

# Set explicit width for column A (closest achievable value to the
# target stored width of 11.7109375 given this engine's column-width
# pixel quantization)
$ws.Columns.Item(1).ColumnWidth = 10.8

# Update the selected range on the sheet
$ws.Range("A8:A13").Select()
